# Enforce HH ID/ExtID format
# Add a new "id_candidate" field to the "model" sheet, right under the
# "Assigned by household form" group (hh_id / hh_country / hh_head), and
# flag it (and any future session-variable rows) via a new
# "isSessionVariable" column. Also makes "model" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Insert a new row 5, copying the formatting (fill/style) of the row above
# (row 4, "hh_head") which belongs to the same "Assigned by household form"
# section, then fill in the new field's name/type.
$ws.Rows("5:5").Insert()
$ws.Range("B4:C4").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)
$ws.Range("B5").Value = "id_candidate"
$ws.Range("C5").Value = "string"

# New column E: isSessionVariable. hh_id (row 2) is the only field flagged.
$ws.Range("E1").Value = "isSessionVariable"
$ws.Range("E5").Value = 1

# Make "model" the active/selected sheet (was "settings").
$ws.Activate()
